$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8579124808311462
$ws.Range("B1").Value = 1.674988985061646
$ws.Range("C1").Value = 6.179018974304199
$ws.Range("D1").Value = 1.887541532516479
$ws.Range("E1").Value = 1.142337441444397
